$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (2-22), columns: D, L, M, N, O, P, Q, R, S, T
$data = @(
    @(44162, "Primera", 100, 4000, 4000, 4000, "`$/bandeja 2 kilos", "Región de O'Higgins", 2000, 2),
    @(44176, "Primera", 150, 3500, 3500, 3500, "`$/bandeja 12 canastillos 125 gramos", "Provincia de Curicó", 2333, 1.5),
    @(44167, "Primera", 500, 3600, 3600, 3600, "`$/bandeja 2 kilos", "Región de O'Higgins", 1800, 2),
    @(44211, "Primera", 40, 2800, 2800, 2800, "`$/bandeja 2 kilos", "Provincia de Linares", 1400, 2),
    @(44211, "Segunda", 30, 2600, 2600, 2600, "`$/bandeja 2 kilos", "Provincia de Linares", 1300, 2),
    @(44210, "Segunda", 150, 2700, 2700, 2700, "`$/bandeja 2 kilos", "Provincia de Linares", 1350, 2),
    @(44200, "Segunda", 50, 2600, 2600, 2600, "`$/bandeja 2 kilos", "Provincia de Linares", 1300, 2),
    @(44265, "Primera", 70, 3600, 3800, 3714, "`$/bandeja 2 kilos", "Provincia de Linares", 1857, 2),
    @(44264, "Primera", 110, 3500, 4000, 3727, "`$/bandeja 2 kilos", "Provincia de Linares", 1864, 2),
    @(44232, "Primera", 60, 3000, 3000, 3000, "`$/bandeja 2 kilos", "Provincia de Linares", 1500, 2),
    @(44169, "Primera", 400, 3600, 3600, 3600, "`$/bandeja 2 kilos", "Provincia de Linares", 1800, 2),
    @(44235, "Primera", 60, 3000, 3000, 3000, "`$/bandeja 2 kilos", "Provincia de Linares", 1500, 2),
    @(44165, "Primera", 400, 3400, 3400, 3400, "`$/bandeja 2 kilos", "Región de O'Higgins", 1700, 2),
    @(44204, "Primera", 50, 3000, 3000, 3000, "`$/bandeja 2 kilos", "Provincia de Linares", 1500, 2),
    @(44204, "Segunda", 140, 2400, 2400, 2400, "`$/bandeja 2 kilos", "Provincia de Linares", 1200, 2),
    @(44186, "Primera", 200, 3000, 3000, 3000, "`$/bandeja 2 kilos", "Provincia de Limarí", 1500, 2),
    @(44202, "Primera", 30, 3000, 3000, 3000, "`$/bandeja 2 kilos", "Provincia de Linares", 1500, 2),
    @(44202, "Segunda", 20, 2600, 2600, 2600, "`$/bandeja 2 kilos", "Provincia de Linares", 1300, 2),
    @(44166, "Primera", 1500, 3600, 3600, 3600, "`$/bandeja 2 kilos", "Región de O'Higgins", 1800, 2),
    @(44172, "Primera", 300, 3400, 3600, 3467, "`$/bandeja 2 kilos", "Provincia de Linares", 1734, 2),
    @(44187, "Primera", 110, 2600, 3000, 2782, "`$/bandeja 2 kilos", "Provincia de Linares", 1391, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 4).Value  = $row[0]   # D - Fecha
    $ws.Cells.Item($r, 12).Value = $row[1]   # L - Calidad
    $ws.Cells.Item($r, 13).Value = $row[2]   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $row[3]   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $row[4]   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $row[5]   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $row[6]   # Q - Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $row[7]   # R - Origen
    $ws.Cells.Item($r, 19).Value = $row[8]   # S - Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $row[9]   # T - Kg / unidad
}
